$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the Price (D) column cells being updated so that
# values like "0.470" / "14.10" / "36.80" keep their exact textual form instead of
# being auto-converted to numbers (which would drop the significant trailing zero).
$dCells = @("D2","D3","D5","D6","D8","D12","D14","D15","D17","D18","D19","D20","D21","D22","D24","D25","D29","D30","D31","D32","D35","D38","D40","D41","D43","D47","D48","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.739.07'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '3.311.73'
$ws.Range("E3").Value = '  +2.26%  '
$ws.Range("D5").Value = '606.83'
$ws.Range("E5").Value = '  +1.84%  '
$ws.Range("D6").Value = '141.84'
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.310.69'
$ws.Range("E8").Value = '  +2.35%  '
$ws.Range("E10").Value = '  +1.72%  '
$ws.Range("E11").Value = '  +2.72%  '
$ws.Range("D12").Value = '0.470'
$ws.Range("E12").Value = '  +1.10%  '
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("D14").Value = '35.04'
$ws.Range("E14").Value = '  +1.97%  '
$ws.Range("D15").Value = '3.857.70'
$ws.Range("E15").Value = '  +2.50%  '
$ws.Range("D17").Value = '3.312.39'
$ws.Range("E17").Value = '  +2.41%  '
$ws.Range("D18").Value = '63.798.14'
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").Value = '6.87'
$ws.Range("E19").Value = '  +1.15%  '
$ws.Range("D20").Value = '481.76'
$ws.Range("E20").Value = '  +1.55%  '
$ws.Range("D21").Value = '14.10'
$ws.Range("E21").Value = '  -0.67%  '
$ws.Range("D22").Value = '0.742'
$ws.Range("E22").Value = '  +1.55%  '
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("D24").Value = '14.00'
$ws.Range("E24").Value = '  +6.14%  '
$ws.Range("D25").Value = '85.30'
$ws.Range("E25").Value = '  +2.06%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = '8.22'
$ws.Range("E29").Value = '  +1.18%  '
$ws.Range("D30").Value = '7.18'
$ws.Range("E30").Value = '  -4.33%  '
$ws.Range("D31").Value = '2.16'
$ws.Range("E31").Value = '  +1.38%  '
$ws.Range("D32").Value = '28.96'
$ws.Range("E32").Value = '  +5.48%  '
$ws.Range("E33").Value = '  -0.79%  '
$ws.Range("E34").Value = '  -0.57%  '
$ws.Range("D35").Value = '1.10'
$ws.Range("E35").Value = '  +1.15%  '
$ws.Range("E36").Value = '  +2.43%  '
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("D38").Value = '0.0₃0748'
$ws.Range("E38").Value = '  +5.43%  '
$ws.Range("E39").Value = '  +2.05%  '
$ws.Range("D40").Value = '435.47'
$ws.Range("E40").Value = '  +3.02%  '
$ws.Range("D41").Value = '3.108.19'
$ws.Range("E41").Value = '  +4.66%  '
$ws.Range("E42").Value = '  +8.18%  '
$ws.Range("D43").Value = '8.35'
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("E45").Value = '  -0.34%  '
$ws.Range("E46").Value = '  +3.16%  '
$ws.Range("D47").Value = '36.80'
$ws.Range("E47").Value = '  +8.72%  '
$ws.Range("D48").Value = '26.40'
$ws.Range("E48").Value = '  +1.71%  '
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("E50").Value = '  -1.00%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '0.114'
$ws.Range("E51").Value = '  -0.44%  '
